$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. Merge the split "Implement those callb" / "acks in index.js and/or dbRequester.js"
#    run (straddled by the _GoBack bookmark) back into a single run with the
#    full sentence, dropping the bookmark from this location (it gets re-added
#    at the very end of the new content below).
$d.Content.Find.Execute(
    "Implement those callbacks in index.js and/or dbRequester.js",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implement those callbacks in index.js and/or dbRequester.js", 2) | Out-Null

# Locate that paragraph again (it is the last paragraph in the document body).
$lastIndex = $d.Paragraphs.Count
$cur = $d.Paragraphs($lastIndex).Range

# NOTE: the engine normalises away an explicit <w:ind> element whenever its
# value equals the value already inherited from the paragraph's numbering
# definition (ListBullet -> left=360 hanging=360), even when it is present in
# the XML passed to InsertXML. To force it to be written out explicitly (as
# the target document requires) we additionally poke
# ParagraphFormat.LeftIndent / FirstLineIndent after inserting the XML; that
# marks the indentation as direct formatting so it gets serialised.

# Paragraph: blank separator line
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr></w:pPr></w:p>")
$cur = $d.Paragraphs($lastIndex).Range
$cur.ParagraphFormat.LeftIndent = 18
$cur.ParagraphFormat.FirstLineIndent = -18

# Paragraph: "19 November 2016 (30mins)" (bold heading)
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>19 November 2016 (30mins)</w:t></w:r></w:p>")
$cur = $d.Paragraphs($lastIndex).Range
$cur.ParagraphFormat.LeftIndent = 18
$cur.ParagraphFormat.FirstLineIndent = -18

# Paragraph: "Implemented call back function successfully. ..."
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/></w:pPr><w:r><w:t>Implemented call back function successfully. Basically this is how it works:</w:t></w:r></w:p>")

# Paragraph: "DB function"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/></w:pPr><w:r><w:t>DB function</w:t></w:r></w:p>")

# Paragraph: "Function RequestMaker (callback) {"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:t xml:space=""preserve"">Function </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>RequestMaker</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> (callback) {</w:t></w:r></w:p>")

# Paragraph: tab + "Make requests to database and put it into 'data' object"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:tab/><w:t>Make requests to database and put it into ‘data’ object</w:t></w:r></w:p>")

# Paragraph: tab + "Callback" + "(data" + ");" + tab + "//This is like a function..."
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:tab/><w:t>Callback</w:t></w:r><w:r><w:t>(data</w:t></w:r><w:r><w:t>);</w:t></w:r><w:r><w:tab/><w:t>//This is like a function that has to be later overloaded by the caller</w:t></w:r></w:p>")

# Paragraph: lastRenderedPageBreak + "}"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>}</w:t></w:r></w:p>")

# Paragraph: "Caller" + tab + tab
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/></w:pPr><w:r><w:t>Caller</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r></w:p>")

# Paragraph: "RequestMaker" (spell-checked) + "(function (data) {"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:proofErr w:type=""spellStart""/><w:r><w:t>RequestMaker</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t>(function (data) {</w:t></w:r></w:p>")

# Paragraph: tab + "Console.log(data);"
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:tab/><w:t>Console.log(data);</w:t></w:r></w:p>")

# Paragraph: "});" followed by the relocated _GoBack bookmark
$cur.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$cur = $d.Paragraphs($lastIndex).Range
$cur.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=""ListBullet""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""0""/></w:numPr><w:ind w:left=""720""/></w:pPr><w:r><w:t>});</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/></w:p>")
